# Updated cryptos list with refreshed Price (D) and Volume(1h) (E) figures.
# Price cells are forced to text via a leading apostrophe (otherwise Excel's
# COM layer auto-coerces strings like "211.01" into floating point numbers,
# corrupting the formatting), then the quote-prefix style introduced by that
# trick is cleared by resetting the cell back to the "Normal" style so the
# cell keeps its original (unstyled) appearance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.669.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "'1.598.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'211.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").Value = "'19.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("D11").Value = "'0.0844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").Value = "'1.591.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "'4.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").Value = "'0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.48%  "
$ws.Range("D16").Value = "'64.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "'26.667.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "'208.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'2.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.68%  "
$ws.Range("D24").Value = "'8.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").Value = "'145.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'7.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.91%  "
$ws.Range("E28").Value = "  +1.79%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").Value = "'0.0506"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").Value = "'0.658"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("D35").Value = "'1.295.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.80%  "
$ws.Range("D36").Value = "'2.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("D39").Value = "'0.844"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.60%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").Value = "'5.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.60%  "
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("D43").Value = "'0.786"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").Value = "'63.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("D45").Value = "'1.736.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("D46").Value = "'0.901"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.41%  "
$ws.Range("D47").Value = "'90.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.21%  "
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").Value = "'7.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.05%  "
